$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.710262656211853
$ws.Range("B1").Value = 1.457847595214844
$ws.Range("C1").Value = 3.973063230514526
$ws.Range("D1").Value = 2.704084157943726
$ws.Range("E1").Value = 0.696803092956543
